$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended after the existing last row (row 47 -> row 48)
$newRow = 48

# Column A holds a date-like string (e.g. "2025-09-29"). Excel's COM layer
# auto-converts such text into a real date serial when assigned directly,
# so force the cell to Text format first to preserve the literal string,
# then restore the default style so no extra formatting is left behind.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-09-29"
$cellA.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "21:22:01"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,633.9357"
